$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").Value = "2024-10-20 00:00:00"
$ws.Range("B79").Value = 73650
$ws.Range("C79").Value = 10320.18
$ws.Range("D79").Value = 9132.91
$ws.Range("E79").Value = 7.1018
